$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be auto-parsed as numbers
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"

# Apply updated cell values
$ws.Range('D2').Value = '70.774.16'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '3.806.25'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '708.56'
$ws.Range('E5').Value = '  +1.96%  '
$ws.Range('D6').Value = '170.49'
$ws.Range('E6').Value = '  -1.87%  '
$ws.Range('D7').Value = '3.806.03'
$ws.Range('E7').Value = '  -1.09%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '0.521'
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('E10').Value = '  -1.59%  '
$ws.Range('D11').Value = '7.44'
$ws.Range('E11').Value = '  +2.46%  '
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('D14').Value = '36.05'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').Value = '4.447.81'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('D16').Value = '3.860.27'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').Value = '70.797.24'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').Value = '7.13'
$ws.Range('E19').Value = '  -1.44%  '
$ws.Range('D20').Value = '17.34'
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').Value = '494.25'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').Value = '10.64'
$ws.Range('E22').Value = '  -4.81%  '
$ws.Range('D23').Value = '0.728'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').Value = '84.41'
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('E25').Value = '  -0.92%  '
$ws.Range('E26').Value = '  -1.81%  '
$ws.Range('E27').Value = '  -1.75%  '
$ws.Range('D28').Value = '3.958.17'
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  -4.81%  '
$ws.Range('E31').Value = '  -3.05%  '
$ws.Range('E32').Value = '  -1.64%  '
$ws.Range('D33').Value = '7.34'
$ws.Range('E33').Value = '  -3.83%  '
$ws.Range('E34').Value = '  -1.95%  '
$ws.Range('D35').Value = '0.172'
$ws.Range('E35').Value = '  -4.20%  '
$ws.Range('D36').Value = '9.14'
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('D37').Value = '3.775.42'
$ws.Range('E37').Value = '  -0.67%  '
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('E39').Value = '  -2.73%  '
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('D41').Value = '2.30'
$ws.Range('E41').Value = '  -3.77%  '
$ws.Range('E43').Value = '  -3.98%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('B46').Value = 'FLOKI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D46').Value = '0.000323'
$ws.Range('E46').Value = '  +5.19%  '
$ws.Range('D47').Value = '164.99'
$ws.Range('E47').Value = '  +0.99%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '48.76'
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = '424.79'
$ws.Range('E49').Value = '  +1.38%  '
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('E51').Value = '  -2.91%  '
